# Auto-generated: updates market-price-derived columns (H-N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1809.8
$ws.Range("I70").Value = 1798
$ws.Range("K70").Value = 5394
$ws.Range("M70").Value = -5124
$ws.Range("H73").Value = 1809.8
$ws.Range("I73").Value = 1798
$ws.Range("K73").Value = 5394
$ws.Range("M73").Value = -4458
$ws.Range("H129").Value = 2098.4546
$ws.Range("J129").Value = 910.75
$ws.Range("L129").Value = 2732.25
$ws.Range("N129").Value = -12732.25
$ws.Range("H132").Value = 19524.64
$ws.Range("I132").Value = 2624.2222
$ws.Range("K132").Value = 7872.6666
$ws.Range("M132").Value = -5342.6666
$ws.Range("H138").Value = 2426.951
$ws.Range("J138").Value = 2972.6829
$ws.Range("L138").Value = 8918.048699999999
$ws.Range("N138").Value = -19198.0487
$ws.Range("H141").Value = 3235.9375
$ws.Range("I141").Value = 1960.8334
$ws.Range("J141").Value = 4001
$ws.Range("K141").Value = 5882.5002
$ws.Range("L141").Value = 12003
$ws.Range("M141").Value = -702.5002000000004
$ws.Range("N141").Value = -22363

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3090.9092
$ws.Range("I2").Value = 2400
$ws.Range("K2").Value = 2400
$ws.Range("M2").Value = -2287
$ws.Range("H32").Value = 11683.547
$ws.Range("I32").Value = 10990.613
$ws.Range("J32").Value = 15071.223
$ws.Range("K32").Value = 10990.613
$ws.Range("L32").Value = 15071.223
$ws.Range("M32").Value = -10703.613
$ws.Range("N32").Value = -15645.223
$ws.Range("H61").Value = 1908.7709
$ws.Range("I61").Value = 1346.3334
$ws.Range("J61").Value = 3146.1333
$ws.Range("K61").Value = 1346.3334
$ws.Range("L61").Value = 3146.1333
$ws.Range("M61").Value = -1134.3334
$ws.Range("N61").Value = -3570.1333
$ws.Range("H74").Value = 1518.9623
$ws.Range("I74").Value = 1243.025
$ws.Range("J74").Value = 2368
$ws.Range("K74").Value = 1243.025
$ws.Range("L74").Value = 2368
$ws.Range("M74").Value = -369.0250000000001
$ws.Range("N74").Value = -4116
$ws.Range("H77").Value = 1518.9623
$ws.Range("I77").Value = 1243.025
$ws.Range("J77").Value = 2368
$ws.Range("K77").Value = 6215.125
$ws.Range("L77").Value = 11840
$ws.Range("M77").Value = -1847.125
$ws.Range("N77").Value = -20576
$ws.Range("H116").Value = 3090.9092
$ws.Range("I116").Value = 2400
$ws.Range("K116").Value = 2400
$ws.Range("M116").Value = -106
$ws.Range("H122").Value = 1795.5294
$ws.Range("I122").Value = 1584.1666
$ws.Range("K122").Value = 4752.4998
$ws.Range("M122").Value = -2302.4998
$ws.Range("H132").Value = 8066474.5
$ws.Range("I132").Value = 11112347
$ws.Range("K132").Value = 33337041
$ws.Range("M132").Value = -33334511
$ws.Range("H136").Value = 1908.7709
$ws.Range("I136").Value = 1346.3334
$ws.Range("J136").Value = 3146.1333
$ws.Range("K136").Value = 4039.0002
$ws.Range("L136").Value = 9438.3999
$ws.Range("M136").Value = -1489.0002
$ws.Range("N136").Value = -14538.3999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3090.9092
$ws.Range("I3").Value = 2400
$ws.Range("K3").Value = 2400
$ws.Range("M3").Value = -2286
$ws.Range("H94").Value = 1071.909
$ws.Range("I94").Value = 1093.6471
$ws.Range("J94").Value = 998
$ws.Range("K94").Value = 1093.6471
$ws.Range("L94").Value = 998
$ws.Range("M94").Value = -642.6470999999999
$ws.Range("N94").Value = -1900
$ws.Range("H134").Value = 1931.36
$ws.Range("I134").Value = 1474.35
$ws.Range("K134").Value = 4423.049999999999
$ws.Range("M134").Value = -1888.049999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3971473.8
$ws.Range("I31").Value = 1488.7
$ws.Range("J31").Value = 6177021
$ws.Range("K31").Value = 1488.7
$ws.Range("L31").Value = 6177021
$ws.Range("M31").Value = -1193.7
$ws.Range("N31").Value = -6177611
$ws.Range("H34").Value = 3971473.8
$ws.Range("I34").Value = 1488.7
$ws.Range("J34").Value = 6177021
$ws.Range("K34").Value = 1488.7
$ws.Range("L34").Value = 6177021
$ws.Range("M34").Value = -1286.7
$ws.Range("N34").Value = -6177425
$ws.Range("H134").Value = 326701.84
$ws.Range("I134").Value = 950.8205
$ws.Range("J134").Value = 3502774.2
$ws.Range("K134").Value = 2852.4615
$ws.Range("L134").Value = 10508322.6
$ws.Range("M134").Value = -317.4615000000003
$ws.Range("N134").Value = -10513392.6

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2664.5312
$ws.Range("I107").Value = 3488.9355
$ws.Range("J107").Value = 2271.3538
$ws.Range("K107").Value = 10466.8065
$ws.Range("L107").Value = 6814.0614
$ws.Range("M107").Value = -8546.806500000001
$ws.Range("N107").Value = -10654.0614

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7815.7896
$ws.Range("J80").Value = 9600
$ws.Range("L80").Value = 9600
$ws.Range("N80").Value = -11596
$ws.Range("H83").Value = 7815.7896
$ws.Range("J83").Value = 9600
$ws.Range("L83").Value = 48000
$ws.Range("N83").Value = -57984
$ws.Range("H126").Value = 26960.916
$ws.Range("I126").Value = 44049.57
$ws.Range("J126").Value = 3036.8
$ws.Range("K126").Value = 132148.71
$ws.Range("L126").Value = 9110.400000000001
$ws.Range("M126").Value = -129678.71
$ws.Range("N126").Value = -14050.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83336890
$ws.Range("I7").Value = 111113600
$ws.Range("J7").Value = 6735
$ws.Range("K7").Value = 111113600
$ws.Range("L7").Value = 6735
$ws.Range("M7").Value = -111113488
$ws.Range("N7").Value = -6959
$ws.Range("H40").Value = 4740.75
$ws.Range("I40").Value = 4219.231
$ws.Range("J40").Value = 7000.6665
$ws.Range("K40").Value = 4219.231
$ws.Range("L40").Value = 7000.6665
$ws.Range("M40").Value = -4083.231
$ws.Range("N40").Value = -7272.6665
$ws.Range("H126").Value = 83336890
$ws.Range("I126").Value = 111113600
$ws.Range("J126").Value = 6735
$ws.Range("K126").Value = 333340800
$ws.Range("L126").Value = 20205
$ws.Range("M126").Value = -333338330
$ws.Range("N126").Value = -25145
$ws.Range("H132").Value = 3969
$ws.Range("I132").Value = 2781.7646
$ws.Range("J132").Value = 5987.3
$ws.Range("K132").Value = 8345.293799999999
$ws.Range("L132").Value = 17961.9
$ws.Range("M132").Value = -5815.293799999999
$ws.Range("N132").Value = -23021.9

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents() | Out-Null
$ws.Range("H126").Value = 4203252
$ws.Range("I126").Value = 4903460.5
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 14710381.5
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -14707911.5
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 1846.2424
$ws.Range("I132").Value = 1241.5
$ws.Range("J132").Value = 2776.6155
$ws.Range("K132").Value = 3724.5
$ws.Range("L132").Value = 8329.8465
$ws.Range("M132").Value = -1194.5
$ws.Range("N132").Value = -13389.8465
$ws.Range("H136").Value = 157070.78
$ws.Range("I136").Value = 204724.33
$ws.Range("J136").Value = 1402.5333
$ws.Range("K136").Value = 614172.99
$ws.Range("L136").Value = 4207.5999
$ws.Range("M136").Value = -611622.99
$ws.Range("N136").Value = -9307.599900000001
